# Set all detector biases to zero in warm condition
# (rows 30-41, columns B:BD on the "Biases" sheet: values go from 1000 -> 0)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Biases")

$ws.Range("B30:BD41").Value = 0

# Reflect the author's final on-screen selection/scroll state:
# topLeftCell moved from A13 to B13, and the selection became
# B30 (active cell) with B30:BD41 selected.
$ws.Range("B30:BD41").Select()
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("B30:BD41").Select()
